# Add a new column J with 2021 data, mirroring the formatting of column I
# (years 2015-2020 already present in columns D-I).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number format, font, borders, alignment, style index)
# from the existing "2020" column (I, rows 4-14) into the new "2021" column (J).
$ws.Range("I4:I14").Copy() | Out-Null
$ws.Range("J4:J14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Populate the new column with the 2021 values.
$ws.Range("J4").Value = 2021
$ws.Range("J5").Value = 1.5
$ws.Range("J6").Value = 0.3
$ws.Range("J7").Value = 0.8
$ws.Range("J8").Value = 0.6
$ws.Range("J9").Value = 1.8
$ws.Range("J10").Value = 0.5
$ws.Range("J11").Value = 0.8
$ws.Range("J12").Value = 1.9
$ws.Range("J13").Value = 4.4000000000000004
$ws.Range("J14").Value = 0.4

# Match the author's recorded selection at the end of the edit.
$ws.Range("L10").Select() | Out-Null
